# Updated main GSC export data:
# The oldest day's row (2025-10-15) is dropped from the rolling "Chart"
# data table; every subsequent row shifts up by one. Deleting the sheet
# row (rather than just clearing values) also drops the now-unused
# "2025-10-15" entry from the shared-strings table and renumbers the
# sheet's used range/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the 2025-10-15 data point (row 1 is the header row).
# EntireRow.Delete shifts rows 3:90 up into 2:89, just like Excel's
# "Delete Sheet Rows" on the row selector.
$ws.Rows.Item(2).EntireRow.Delete()
